$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Update Details" text first, then the date, so that the
# shared-string table regenerates them in the same order as the target
# workbook (long text before the date).
$ws.Range("F3").Value = "Initial setup complete (Create, Edit, Delete, Details)." + [char]10 + "Firebase Connected." + [char]10 + "Model Class Finalised." + [char]10 + "Authentication Added." + [char]10 + "Model Class Automatically Collects Totals."
$ws.Range("E3").Value = "22/11/2021"

# The extra two lines of update text make row 3 noticeably taller.
$ws.Rows(3).RowHeight = 90.75

# Move/keep the active selection on the worksheet.
$ws.Range("F14").Select()
